$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45042
$ws.Range("J2").Value = 60
$ws.Range("N2").Value = '$/caja 18 kilos empedrada'
# Row 3
$ws.Range("D3").Value = 44280
$ws.Range("J3").Value = 30
$ws.Range("N3").Value = '$/caja 18 kilos empedrada'
$ws.Range("P3").Value = 1389
$ws.Range("Q3").Value = 18
# Row 5
$ws.Range("D5").Value = 45041
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 24000
$ws.Range("L5").Value = 24000
$ws.Range("M5").Value = 24000
$ws.Range("P5").Value = 1333
# Row 6
$ws.Range("D6").Value = 45015
$ws.Range("J6").Value = 50
# Row 7
$ws.Range("D7").Value = 45033
$ws.Range("J7").Value = 80
$ws.Range("N7").Value = '$/caja 18 kilos granel'
# Row 8
$ws.Range("D8").Value = 45043
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 24000
$ws.Range("L8").Value = 24000
$ws.Range("M8").Value = 24000
$ws.Range("N8").Value = '$/caja 18 kilos empedrada'
$ws.Range("P8").Value = 1333
$ws.Range("Q8").Value = 18
# Row 9
$ws.Range("D9").Value = 45014
$ws.Range("J9").Value = 30
$ws.Range("K9").Value = 8000
$ws.Range("L9").Value = 8000
$ws.Range("M9").Value = 8000
$ws.Range("N9").Value = '$/caja 18 kilos empedrada'
$ws.Range("P9").Value = 444
$ws.Range("Q9").Value = 18
# Row 10
$ws.Range("D10").Value = 44285
$ws.Range("J10").Value = 20
$ws.Range("K10").Value = 25000
$ws.Range("L10").Value = 25000
$ws.Range("M10").Value = 25000
$ws.Range("P10").Value = 1389
# Row 11
$ws.Range("D11").Value = 44315
$ws.Range("I11").Value = 'Especial'
$ws.Range("J11").Value = 10
$ws.Range("K11").Value = 30000
$ws.Range("L11").Value = 30000
$ws.Range("M11").Value = 30000
$ws.Range("N11").Value = '$/caja 20 kilos empedrada'
$ws.Range("P11").Value = 1500
$ws.Range("Q11").Value = 20
# Row 12
$ws.Range("D12").Value = 44315
$ws.Range("J12").Value = 20
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 15000
$ws.Range("N12").Value = '$/caja 15 kilos granel'
$ws.Range("P12").Value = 1000
$ws.Range("Q12").Value = 15
# Row 13
$ws.Range("D13").Value = 44313
$ws.Range("J13").Value = 20
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 15000
$ws.Range("N13").Value = '$/caja 15 kilos empedrada'
$ws.Range("Q13").Value = 15
# Row 14
$ws.Range("D14").Value = 44313
$ws.Range("J14").Value = 20
$ws.Range("K14").Value = 30000
$ws.Range("L14").Value = 30000
$ws.Range("M14").Value = 30000
$ws.Range("N14").Value = '$/caja 20 kilos empedrada'
$ws.Range("P14").Value = 1500
$ws.Range("Q14").Value = 20
# Row 15
$ws.Range("D15").Value = 45040
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 18000
$ws.Range("L15").Value = 18000
$ws.Range("M15").Value = 18000
$ws.Range("P15").Value = 1000
# Row 16
$ws.Range("D16").Value = 45034
$ws.Range("N16").Value = '$/caja 18 kilos granel'
# Row 17
$ws.Range("D17").Value = 45044
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 40
$ws.Range("K17").Value = 24000
$ws.Range("L17").Value = 24000
$ws.Range("M17").Value = 24000
$ws.Range("N17").Value = '$/caja 18 kilos empedrada'
$ws.Range("P17").Value = 1333
$ws.Range("Q17").Value = 18
# Row 18
$ws.Range("D18").Value = 44293
$ws.Range("J18").Value = 10
$ws.Range("K18").Value = 25000
$ws.Range("L18").Value = 25000
$ws.Range("M18").Value = 25000
$ws.Range("N18").Value = '$/caja 15 kilos empedrada'
$ws.Range("P18").Value = 1667

